# Weekly fruit/vegetable price update for "Feria Lagunitas de Puerto Montt - Mango".
# A new price record (week of 2022-07-05) is inserted as row 206, pushing the
# existing rows 206:212 down to 207:213 (dimension grows from T212 to T213).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 206, shifting rows 206:212 down to 207:213.
$ws.Rows("206").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A206").Value = 4
$ws.Range("B206").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C206").Value = "Los Lagos"
$ws.Range("D206").Value = 44747
$ws.Range("E206").Value = 10
$ws.Range("F206").Value = "Fruta"
$ws.Range("G206").Value = 100108
$ws.Range("H206").Value = "Tropicales y subtropicales"
$ws.Range("I206").Value = 100108002
$ws.Range("J206").Value = "Mango"
$ws.Range("K206").Value = "Sin especificar"
$ws.Range("L206").Value = "Primera"
$ws.Range("M206").Value = 200
$ws.Range("N206").Value = 9500
$ws.Range("O206").Value = 9500
$ws.Range("P206").Value = 9500
$ws.Range("Q206").Value = "$/bandeja 4 kilos"
$ws.Range("R206").Value = "Perú"
$ws.Range("S206").Value = 2375
$ws.Range("T206").Value = 4
